# Weekly update: insert a new data row at the top of the price list block
# (row 150), shifting the existing rows 150:174 down to 151:175, then fill
# in the new row with this week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 150; Excel shifts rows 150:174 down to 151:175
# and copies the formatting (including the date number format on column D)
# from the row above.
$ws.Rows("150:150").Insert()

# Populate the new row 150 with the latest week's record.
$ws.Range("A150").Value = 1
$ws.Range("B150").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C150").Value = "Arica y Parinacota"
$ws.Range("D150").Value = 45131
$ws.Range("E150").Value = 15
$ws.Range("F150").Value = 100112042
$ws.Range("G150").Value = "Locoto"
$ws.Range("H150").Value = "Sin especificar"
$ws.Range("I150").Value = "Primera"
$ws.Range("J150").Value = 200
$ws.Range("K150").Value = 30000
$ws.Range("L150").Value = 32000
$ws.Range("M150").Value = 31000
$ws.Range("N150").Value = "$/caja 20 kilos"
$ws.Range("O150").Value = "Región de Arica y Parinacota"
$ws.Range("P150").Value = 1550
$ws.Range("Q150").Value = 20
$ws.Range("R150").Value = "Hortaliza"
